$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.949.42"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "3.244.74"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "395.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.60%  "

$ws.Range("E7").Value = "  +3.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -1.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.30%  "

$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("D13").Value = "3.754.09"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "

$ws.Range("D16").Value = "3.251.84"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.49%  "

$ws.Range("D19").Value = "56.734.22"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("E20").Value = "  -3.58%  "

$ws.Range("E21").Value = "  +4.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "291.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.81%  "

$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.53%  "

$ws.Range("E32").Value = "  -2.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("E34").Value = "  +10.61%  "

$ws.Range("E35").Value = "  -2.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("E39").Value = "  -3.47%  "

$ws.Range("E40").Value = "  -4.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "137.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.06%  "

$ws.Range("E42").Value = "  +1.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.284"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.12%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "

$ws.Range("D49").Value = "2.146.39"
$ws.Range("E49").Value = "  -0.78%  "

$ws.Range("E50").Value = "  -4.38%  "

$ws.Range("E51").Value = "  -8.14%  "
